$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 86600
$ws.Range("J3").Value = 86600
$ws.Range("L3").Value = 86600
$ws.Range("N3").Value = -86828
$ws.Range("H88").Value = 7575.6665
$ws.Range("I88").Value = 7600
$ws.Range("J88").Value = 7454
$ws.Range("K88").Value = 7600
$ws.Range("L88").Value = 7454
$ws.Range("M88").Value = -7194
$ws.Range("N88").Value = -8266
$ws.Range("H91").Value = 7575.6665
$ws.Range("I91").Value = 7600
$ws.Range("J91").Value = 7454
$ws.Range("K91").Value = 7600
$ws.Range("L91").Value = 7454
$ws.Range("M91").Value = -6196
$ws.Range("N91").Value = -10262
$ws.Range("H98").Value = 942.8333
$ws.Range("I98").Value = 1026.5714
$ws.Range("K98").Value = 1026.5714
$ws.Range("M98").Value = 471.4286
$ws.Range("H100").Value = 5450.174
$ws.Range("I100").Value = 5207.3
$ws.Range("J100").Value = 5637
$ws.Range("K100").Value = 5207.3
$ws.Range("L100").Value = 5637
$ws.Range("M100").Value = -4666.3
$ws.Range("N100").Value = -6719
$ws.Range("H102").Value = 86600
$ws.Range("J102").Value = 86600
$ws.Range("L102").Value = 86600
$ws.Range("N102").Value = -93090
$ws.Range("H111").Value = 2700.2
$ws.Range("I111").Value = 1730.7142
$ws.Range("J111").Value = 4962.3335
$ws.Range("K111").Value = 5192.142599999999
$ws.Range("L111").Value = 14887.0005
$ws.Range("M111").Value = -2125.142599999999
$ws.Range("N111").Value = -21021.0005
$ws.Range("H122").Value = 942.8333
$ws.Range("I122").Value = 1026.5714
$ws.Range("K122").Value = 3079.7142
$ws.Range("M122").Value = -629.7142000000003
$ws.Range("H132").Value = 15284.069
$ws.Range("I132").Value = 2480.5293
$ws.Range("K132").Value = 7441.5879
$ws.Range("M132").Value = -4911.5879
$ws.Range("H135").Value = 1769.3077
$ws.Range("I135").Value = 1647.7567
$ws.Range("K135").Value = 14829.8103
$ws.Range("M135").Value = -12294.8103
$ws.Range("H138").Value = 3554.7722
$ws.Range("J138").Value = 4067.15
$ws.Range("L138").Value = 12201.45
$ws.Range("N138").Value = -22481.45
$ws.Range("H141").Value = 5936.6924
$ws.Range("I141").Value = 2897.9092
$ws.Range("K141").Value = 8693.7276
$ws.Range("M141").Value = -3513.7276

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18187076
$ws.Range("I32").Value = 18187076
$ws.Range("K32").Value = 18187076
$ws.Range("M32").Value = -18186789
$ws.Range("H63").Value = 8760.875
$ws.Range("I63").Value = 5698.8
$ws.Range("K63").Value = 5698.8
$ws.Range("M63").Value = -5012.8
$ws.Range("H66").Value = 8760.875
$ws.Range("I66").Value = 5698.8
$ws.Range("K66").Value = 28494
$ws.Range("M66").Value = -25062
$ws.Range("H97").Value = 2538.875
$ws.Range("J97").Value = 5241.6
$ws.Range("L97").Value = 5241.6
$ws.Range("N97").Value = -6233.6
$ws.Range("H106").Value = 54956.668
$ws.Range("J106").Value = 54956.668
$ws.Range("L106").Value = 54956.668
$ws.Range("N106").Value = -57480.668
$ws.Range("H132").Value = 2619.0417
$ws.Range("I132").Value = 2564.238
$ws.Range("K132").Value = 7692.714
$ws.Range("M132").Value = -5162.714

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 23000
$ws.Range("J103").Value = 23000
$ws.Range("L103").Value = 23000
$ws.Range("N103").Value = -25344

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1803.5405
$ws.Range("I31").Value = 1829.9032
$ws.Range("J31").Value = 1667.3334
$ws.Range("K31").Value = 1829.9032
$ws.Range("L31").Value = 1667.3334
$ws.Range("M31").Value = -1534.9032
$ws.Range("N31").Value = -2257.3334
$ws.Range("H34").Value = 1803.5405
$ws.Range("I34").Value = 1829.9032
$ws.Range("J34").Value = 1667.3334
$ws.Range("K34").Value = 1829.9032
$ws.Range("L34").Value = 1667.3334
$ws.Range("M34").Value = -1627.9032
$ws.Range("N34").Value = -2071.3334
$ws.Range("H43").Value = 38600
$ws.Range("J43").Value = 38600
$ws.Range("L43").Value = 38600
$ws.Range("N43").Value = -38968
$ws.Range("H94").Value = 10630.363
$ws.Range("I94").Value = 20649.6
$ws.Range("K94").Value = 20649.6
$ws.Range("M94").Value = -20198.6
$ws.Range("H101").Value = 38600
$ws.Range("J101").Value = 38600
$ws.Range("L101").Value = 38600
$ws.Range("N101").Value = -45090
$ws.Range("H132").Value = 2052.923
$ws.Range("I132").Value = 2052.923
$ws.Range("K132").Value = 6158.768999999999
$ws.Range("M132").Value = -3628.768999999999
$ws.Range("H138").Value = 66702.08
$ws.Range("J138").Value = 70000
$ws.Range("L138").Value = 70000
$ws.Range("N138").Value = -80280

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 85.36364
$ws.Range("I2").Value = 88
$ws.Range("K2").Value = 528
$ws.Range("M2").Value = -415
$ws.Range("H23").Value = 1625.6
$ws.Range("I23").Value = 3881.5
$ws.Range("J23").Value = 121.666664
$ws.Range("K23").Value = 11644.5
$ws.Range("L23").Value = 364.999992
$ws.Range("M23").Value = -11409.5
$ws.Range("N23").Value = -834.999992
$ws.Range("H62").Value = 15000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 15000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H92").Value = 575.1111
$ws.Range("I92").Value = 259
$ws.Range("J92").Value = 733.1667
$ws.Range("K92").Value = 777
$ws.Range("L92").Value = 2199.5001
$ws.Range("M92").Value = 471
$ws.Range("N92").Value = -4695.5001
$ws.Range("H107").Value = 3364.2856
$ws.Range("J107").Value = 4600
$ws.Range("L107").Value = 13800
$ws.Range("N107").Value = -17640
$ws.Range("H121").Value = 2068.5
$ws.Range("J121").Value = 2363.9
$ws.Range("L121").Value = 7091.700000000001
$ws.Range("N121").Value = -9711.700000000001
$ws.Range("H122").Value = 757.46155
$ws.Range("I122").Value = 871.44446
$ws.Range("J122").Value = 501
$ws.Range("K122").Value = 7843.00014
$ws.Range("L122").Value = 4509
$ws.Range("M122").Value = -5393.00014
$ws.Range("N122").Value = -9409
$ws.Range("H137").Value = 2880.8696
$ws.Range("I137").Value = 2543.3635
$ws.Range("J137").Value = 3190.25
$ws.Range("K137").Value = 7630.0905
$ws.Range("L137").Value = 9570.75
$ws.Range("M137").Value = -2530.0905
$ws.Range("N137").Value = -19770.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4536.625
$ws.Range("I70").Value = 4465.6665
$ws.Range("K70").Value = 4465.6665
$ws.Range("M70").Value = -4195.6665
$ws.Range("H73").Value = 4536.625
$ws.Range("I73").Value = 4465.6665
$ws.Range("K73").Value = 4465.6665
$ws.Range("M73").Value = -3529.6665
$ws.Range("H80").Value = 7615.5864
$ws.Range("I80").Value = 6473.1577
$ws.Range("J80").Value = 9786.200000000001
$ws.Range("K80").Value = 6473.1577
$ws.Range("L80").Value = 9786.200000000001
$ws.Range("M80").Value = -5475.1577
$ws.Range("N80").Value = -11782.2
$ws.Range("H83").Value = 7615.5864
$ws.Range("I83").Value = 6473.1577
$ws.Range("J83").Value = 9786.200000000001
$ws.Range("K83").Value = 32365.7885
$ws.Range("L83").Value = 48931
$ws.Range("M83").Value = -27373.7885
$ws.Range("N83").Value = -58915
$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4209.6665
$ws.Range("I7").Value = 3222.5557
$ws.Range("K7").Value = 3222.5557
$ws.Range("M7").Value = -3110.5557
$ws.Range("H46").Value = 4077.4443
$ws.Range("J46").Value = 4499.625
$ws.Range("L46").Value = 4499.625
$ws.Range("N46").Value = -4875.625
$ws.Range("H126").Value = 4209.6665
$ws.Range("I126").Value = 3222.5557
$ws.Range("K126").Value = 9667.667099999999
$ws.Range("M126").Value = -7197.667099999999
$ws.Range("H132").Value = 8402
$ws.Range("I132").Value = 2869.3333
$ws.Range("J132").Value = 25000
$ws.Range("K132").Value = 8607.999899999999
$ws.Range("L132").Value = 75000
$ws.Range("M132").Value = -6077.999899999999
$ws.Range("N132").Value = -80060
$ws.Range("H138").Value = 69996.664
$ws.Range("J138").Value = 69996.664
$ws.Range("L138").Value = 69996.664
$ws.Range("N138").Value = -80276.664

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2178.0952
$ws.Range("I136").Value = 1127.5
$ws.Range("K136").Value = 3382.5
$ws.Range("M136").Value = -832.5
